# Add a new worksheet "No Registrados" at the end of the workbook, listing
# attendees who joined the call but were not registered in the course roster.
$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet, so it lands at
# the end of the tab strip (Worksheets.Add() alone would insert at the front).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "No Registrados"

# Header row
$headers = @('Nombre','Apellido','Correo electrónico','Hora a la que se unió','Hora a la que salió','NombreCompleto','duración')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Copy the header formatting (bold, centered, thin border) used by the other
# attendance sheets so the new sheet matches the existing look.
$srcHeader = $wb.Worksheets.Item("Profesionales").Range("A1:G1")
$srcHeader.Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows
$data = @(
    ,@('Javier','Ballesteros','javier.ballesteros@uptc.edu.co','2023-09-06 10:02:00','2023-09-06 10:30:00','Javier Ballesteros',28)
    ,@('Mónica','Quevedo','moni**********@***.com','2023-09-06 09:46:00','2023-09-06 11:29:00','Mónica Quevedo',103)
    ,@('Cristian David','Saenz Gamboa','cristian.saenz01@uptc.edu.co','2023-09-06 10:05:00','2023-09-06 12:02:00','Cristian David Saenz Gamboa',117)
    ,@('Telematics','Uptc','telematics@uptc.edu.co','2023-09-06 10:03:00','2023-09-06 12:08:00','Telematics Uptc',125)
    ,@('Brayan','Borda','brayan.borda01@uptc.edu.co','2023-09-06 10:02:00','2023-09-06 12:04:00','Brayan Borda01',122)
    ,@('Lorena','Larotta','tele*****************@***.com','2023-09-06 10:04:00','2023-09-06 10:40:00','Lorena Larotta',36)
    ,@('Santigo','Rodriguez Villarruel','sant*****************@***.com','2023-09-06 10:02:00','2023-09-06 11:59:00','Santigo Rodriguez Villarruel',117)
    ,@('Norbey','Rojas','norbey.rojas@uptc.edu.co','2023-09-06 10:03:00','2023-09-06 12:04:00','Norbey Rojas',121)
    ,@('Cristian David','Saenz Gamboa','cristian.saenz01@uptc.edu.co','2023-09-06 10:17:00','2023-09-06 11:58:00','Cristian David Saenz Gamboa',101)
    ,@('Telematics','Uptc','telematics@uptc.edu.co','2023-09-06 10:03:00','2023-09-06 12:04:00','Telematics Uptc',121)
    ,@('  ','Nan','','2023-09-06 10:02:00','2023-09-06 11:58:00','+57 311 *****26',116)
    ,@('Adquisiciones','Biblioteca','adquisiciones.biblioteca@uptc.edu.co','2023-09-06 09:57:00','2023-09-06 09:57:00','Adquisiciones Biblioteca',0)
    ,@('Filyp','Garoch','inge********@***.com','2023-09-06 10:03:00','2023-09-06 11:47:00','Filyp Garoch',104)
    ,@('Mónica','Quevedo','moni**********@***.com','2023-09-06 10:03:00','2023-09-06 11:53:00','Mónica Quevedo',110)
    ,@('Norbey','Rojas','norbey.rojas@uptc.edu.co','2023-09-06 10:09:00','2023-09-06 11:53:00','Norbey Rojas',104)
    ,@('Cristian David','Saenz Gamboa','cristian.saenz01@uptc.edu.co','2023-09-06 10:09:00','2023-09-06 11:39:00','Cristian David Saenz Gamboa',90)
    ,@('Filyp','Garoch','inge********@***.com','2023-09-06 10:10:00','2023-09-06 11:11:00','Filyp Garoch',61)
    ,@('Mónica','Quevedo','moni**********@***.com','2023-09-06 09:52:00','2023-09-06 11:12:00','Mónica Quevedo',80)
    ,@('Santigo','Rodriguez Villarruel','sant*****************@***.com','2023-09-06 10:07:00','2023-09-06 11:12:00','Santigo Rodriguez Villarruel',65)
    ,@('Norbey','Rojas','norbey.rojas@uptc.edu.co','2023-09-06 10:59:00','2023-09-06 11:12:00','Norbey Rojas',13)
    ,@('Cristian David','Saenz Gamboa','cristian.saenz01@uptc.edu.co','2023-09-06 10:08:00','2023-09-06 11:12:00','Cristian David Saenz Gamboa',64)
    ,@('Sandra','García','sand*********@***.com','2023-09-06 09:59:00','2023-09-06 11:56:00','Sandra García',117)
    ,@('Orlando','Medina Castro','orlando.medina@uptc.edu.co','2023-09-06 10:04:00','2023-09-06 11:57:00','Orlando Medina Castro',113)
    ,@('Cristian David','Saenz Gamboa','cristian.saenz01@uptc.edu.co','2023-09-06 10:09:00','2023-09-06 11:17:00','Cristian David Saenz Gamboa',68)
    ,@('Alejandro','Enciso','aenc***@***.com','2023-09-06 09:58:00','2023-09-06 11:55:00','Alejandro Enciso',117)
    ,@('Orlando','Medina Castro','orlando.medina@uptc.edu.co','2023-09-06 10:10:00','2023-09-06 11:56:00','Orlando Medina Castro',106)
    ,@('Norbey','Rojas','norbey.rojas@uptc.edu.co','2023-09-06 10:57:00','2023-09-06 12:00:00','Norbey Rojas',63)
    ,@('Cristian David','Saenz Gamboa','cristian.saenz01@uptc.edu.co','2023-09-06 10:31:00','2023-09-06 11:52:00','Cristian David Saenz Gamboa',81)
    ,@('  ','Nan','','2023-09-06 10:20:00','2023-09-06 10:21:00','+57 311 *****26',1)
    ,@('Orlando','Medina Castro','orlando.medina@uptc.edu.co','2023-09-06 10:24:00','2023-09-06 11:47:00','Orlando Medina Castro',83)
    ,@('Mónica','Quevedo','moni**********@***.com','2023-09-06 10:30:00','2023-09-06 11:47:00','Mónica Quevedo',77)
    ,@('Monica','Quevedo Reyes','mque******@***.com','2023-09-06 09:54:00','2023-09-06 10:25:00','Monica Quevedo Reyes',31)
    ,@('Sandra','García','sand*********@***.com','2023-09-06 10:18:00','2023-09-06 11:51:00','Sandra García',93)
    ,@('Orlando','Medina Castro','orlando.medina@uptc.edu.co','2023-09-06 09:58:00','2023-09-06 11:51:00','Orlando Medina Castro',113)
    ,@('Mónica','Quevedo','moni**********@***.com','2023-09-06 09:59:00','2023-09-06 11:51:00','Mónica Quevedo',112)
    ,@('Edgar David','Giraldo Valencia','edga*********@***.com','2023-09-06 09:55:00','2023-09-06 11:26:00','Edgar David Giraldo Valencia',91)
    ,@('Orlando','Medina Castro','orlando.medina@uptc.edu.co','2023-09-06 10:09:00','2023-09-06 11:11:00','Orlando Medina Castro',62)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
